$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8: tighten totals, add "Expectation" note ---
$ws.Range("H8").Value = 25000
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = "Cull count will be >0 but will not increase as much as the above"

# --- Row 7: fix total ---
$ws.Range("K7").Value = 25000

# --- Row 9: "Test scenario" -> "Sine wave", new starting value, fixed total ---
$ws.Range("B9").Value = "Sine wave"
$ws.Range("E9").Value = 497.7
$ws.Range("H9").Value = 25000
$ws.Range("K9").Value = 25000

# --- Row 10: "Test scenario" -> "Sine wave", new description, fixed total, new expectation ---
$ws.Range("B10").Value = "Sine wave"
$ws.Range("D10").Value = "As above, but the wavelength is shorter (higher frequency) and the amplitude is smaller"
$ws.Range("H10").Value = 25000
$ws.Range("K10").Value = 25000
$ws.Range("L10").Value = "Similar number of trees as above"

# --- Row 11: "Test scenario" -> "Random complex wave", fill manager-budget columns, update description, fixed total ---
$ws.Range("B11").Value = "Random complex wave"
$ws.Range("D11").Value = "Manager budget increases and decreases unpredictably (using Fourier series to create random complex waves), user budget increases linearly. There are currently 10 waves that constitute this scenario"
$ws.Range("E11").Value = "449.9 - 527.7"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "variable"
$ws.Range("H11").Value = 25000
$ws.Range("K11").Value = 25000

# --- Remove the old "T3 / Manager budget increases in steps..." scenario row entirely ---
$ws.Rows("12:12").Delete()

# The "Not sure I want to include the above orange scenario" note (now on row 12
# after the delete above) is dropped, but the row/cell formatting stays.
$ws.Range("A12").ClearContents()

# --- Remove the stray long-form comment row (now row 15 after the shift above) ---
$ws.Rows("15:15").Delete()

# --- Update the manager/user budget ratio notes ---
$ws.Range("A17").Value = "manager budget not more than 140% of user budget"
$ws.Range("A18").Value = "manager budget no less than 70% of user budget"

# Restore the reported selection
$ws.Range("D11").Select()
